# Build site at 2023-04-12 14:53:07 UTC
# Fill in missing / corrected content on the LOB1262 syllabus sheet and
# insert two new rows (Docentes responsaveis entries), shifting several
# rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos (PT) text was missing - it previously held the
#     professor name by mistake. Fix it with the real objectives text.
$ws.Cells.Item(10, 2).Value = "Propiciar ao discente conhecimento dos fundamentos da Educação Ambiental utilizando como base os problemas ambientais da atualidade. Desenvolver atividades práticas integradas à região. Orientar o desenvolvimento de projetos relacionados à Gestão e Educação Ambiental."
$ws.Cells.Item(10, 3).Value = "Propiciar ao discente conhecimento dos fundamentos da Educação Ambiental utilizando como base os problemas ambientais da atualidade. Desenvolver atividades práticas integradas à região. Orientar o desenvolvimento de projetos relacionados à Gestão e Educação Ambiental."

# --- Insert two new rows after row 12 ("Docentes responsáveis:") to hold
#     the two professor names (previously misplaced further down). This
#     shifts every following row down by 2 (13->15 .. 21->23), which is
#     exactly the row layout of the target sheet.
$ws.Rows.Item(13).Resize(2).Insert()

# The inserted rows inherit row 12's (bold, label) formatting - reset the
# B/C cell formatting to the regular body style used elsewhere (copy it
# from row 11, which already has the plain wrap-text/red-text styles).
$ws.Range("B11:C11").Copy()
$ws.Range("B13:C14").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# New row 13: first professor, no label in column A.
$ws.Cells.Item(13, 1).Clear()
$ws.Cells.Item(13, 2).Value = "9146830 - Danúbia Caporusso Bargos"
$ws.Cells.Item(13, 3).Value = "9146830 - Danúbia Caporusso Bargos"

# New row 14: second professor, no label in column A.
$ws.Cells.Item(14, 1).Clear()
$ws.Cells.Item(14, 2).Value = "5817650 - Érica Leonor Romão"
$ws.Cells.Item(14, 3).Value = "5817650 - Érica Leonor Romão"

# --- Row 15 ("Programa resumido:") was missing its PT summary text
#     (had the activation date by mistake).
$ws.Cells.Item(15, 2).Value = "Considerações gerais sobre a problemática ambiental. Evolução das questões ambientais no Brasil e no mundo. Educação e Gestão Ambiental. Elaboração e acompanhamento de projetos de educação ambiental."
$ws.Cells.Item(15, 3).Value = "Considerações gerais sobre a problemática ambiental. Evolução das questões ambientais no Brasil e no mundo. Educação e Gestão Ambiental. Elaboração e acompanhamento de projetos de educação ambiental."

# --- Row 17 ("Programa:") was missing its PT syllabus text (had the
#     professor name by mistake).
$ws.Cells.Item(17, 2).Value = "Sociedade, natureza e desenvolvimento. A relação degradação ambiental-qualidade de vida. Meio ambiente e cidadania. Percepção e Interpretação ambiental. Meio ambiente e representação social. Histórico da educação ambiental e conceitos de meio ambiente; Conceitos, princípios e pensamentos norteadores da Educação Ambiental. A questão ambiental e as conferências mundiais de meio ambiente. O movimento ambientalista e o histórico da EA no Brasil e no mundo; A Agenda 21 e educação ambiental. A política nacional de educação ambiental (pnea) e legislação correlata: A abordagem interdisciplinar da educação ambiental; Educação como instrumento de Gestão Ambiental. Educação ambiental nas empresas e o Sistema de Gestão Ambiental. Projetos, reflexões e práticas da Educação Ambiental. Análise e vivências de experiências práticas de educação ambiental em diferentes contextos. Metodologia de projetos, oficinas e capacitação em educação ambiental."
$ws.Cells.Item(17, 3).Value = "Sociedade, natureza e desenvolvimento. A relação degradação ambiental-qualidade de vida. Meio ambiente e cidadania. Percepção e Interpretação ambiental. Meio ambiente e representação social. Histórico da educação ambiental e conceitos de meio ambiente; Conceitos, princípios e pensamentos norteadores da Educação Ambiental. A questão ambiental e as conferências mundiais de meio ambiente. O movimento ambientalista e o histórico da EA no Brasil e no mundo; A Agenda 21 e educação ambiental. A política nacional de educação ambiental (pnea) e legislação correlata: A abordagem interdisciplinar da educação ambiental; Educação como instrumento de Gestão Ambiental. Educação ambiental nas empresas e o Sistema de Gestão Ambiental. Projetos, reflexões e práticas da Educação Ambiental. Análise e vivências de experiências práticas de educação ambiental em diferentes contextos. Metodologia de projetos, oficinas e capacitação em educação ambiental."

# --- Row 20 ("Método:") was missing its evaluation text (had a
#     professor name by mistake).
$ws.Cells.Item(20, 2).Value = "Avaliação baseada em provas, exercícios, projetos, seminários e outras formas de avaliação, sendo a nota final correspondente a média ponderada das notas atribuídas às avaliações aplicadas"
$ws.Cells.Item(20, 3).Value = "Avaliação baseada em provas, exercícios, projetos, seminários e outras formas de avaliação, sendo a nota final correspondente a média ponderada das notas atribuídas às avaliações aplicadas"

# --- Row 21 ("Critério:") was missing its pass criteria text.
$ws.Cells.Item(21, 2).Value = "Nota Final: NF ≥ 5,0"
$ws.Cells.Item(21, 3).Value = "Nota Final: NF ≥ 5,0"

# --- Row 22 ("Norma de recuperação:") was missing its make-up exam text.
$ws.Cells.Item(22, 2).Value = "Provas e/ou exercícios dirigidos"
$ws.Cells.Item(22, 3).Value = "Provas e/ou exercícios dirigidos"

# --- Row 23 ("Bibliografia:") was missing its bibliography text entirely.
$ws.Cells.Item(23, 2).Value = "CARVALHO, I. C. M.; Educação Ambiental e formação do sujeito ecológico. São Paulo: Cortez, 2006.CINQUETTI, H. C. S.; LOGAREZZI, A. (Org.). Consumo e Resíduo - Fundamentos para o trabalho educativo. 1 ed. São Carlos: EdUFSCar, 2006, v. 1.DIAS, G. F. Dinâmica e instrumentação para educação ambiental. 1. ed. São Paulo: Gaia, 2010. v. 1. 216p.DIAS, G. F. Educação e Gestão Ambiental. 1. ed. São Paulo: Editora Gaia Ltda, 2006. v. 1. 118p.DIAS, G. F. Educação Ambiental: princípios e práticas. 6a ed. São Paulo: Gaia, 2000.GUIMARÃES, M. (org.) Caminhos da educação ambiental: da forma à ação. Campinas, SP: Papirus, 2006.JACOBI, Pedro Roberto, MONTEIRO,F. M ; FERNANDES, M. L. B. . Educação e Sustentabilidade- caminhos e práticas para uma educação transformadora. São Paulo: Evoluir Cultural, 2009. v. 01. 108p.JACOBI, Pedro Roberto OLIVEIRA, F. C. J. F. (Org.). Educação, Meio Ambiente e Cidadania - reflexões e experiências. São Paulo: SMA/CEAM, 1998. 121p LOUREIRO, C. F. B. Trajetória e fundamentos da educação ambiental. 4. ed. São Paulo: Cortez editora, 2012. 165pPHILIPPI JR., A & PELICIONI, M. C. F. (Eds). 2005. Educação ambiental e sustentabilidade. Barueri SP: Manole. 878p. (Coleção Ambiental, 3)."
$ws.Cells.Item(23, 3).Value = "CARVALHO, I. C. M.; Educação Ambiental e formação do sujeito ecológico. São Paulo: Cortez, 2006.CINQUETTI, H. C. S.; LOGAREZZI, A. (Org.). Consumo e Resíduo - Fundamentos para o trabalho educativo. 1 ed. São Carlos: EdUFSCar, 2006, v. 1.DIAS, G. F. Dinâmica e instrumentação para educação ambiental. 1. ed. São Paulo: Gaia, 2010. v. 1. 216p.DIAS, G. F. Educação e Gestão Ambiental. 1. ed. São Paulo: Editora Gaia Ltda, 2006. v. 1. 118p.DIAS, G. F. Educação Ambiental: princípios e práticas. 6a ed. São Paulo: Gaia, 2000.GUIMARÃES, M. (org.) Caminhos da educação ambiental: da forma à ação. Campinas, SP: Papirus, 2006.JACOBI, Pedro Roberto, MONTEIRO,F. M ; FERNANDES, M. L. B. . Educação e Sustentabilidade- caminhos e práticas para uma educação transformadora. São Paulo: Evoluir Cultural, 2009. v. 01. 108p.JACOBI, Pedro Roberto OLIVEIRA, F. C. J. F. (Org.). Educação, Meio Ambiente e Cidadania - reflexões e experiências. São Paulo: SMA/CEAM, 1998. 121p LOUREIRO, C. F. B. Trajetória e fundamentos da educação ambiental. 4. ed. São Paulo: Cortez editora, 2012. 165pPHILIPPI JR., A & PELICIONI, M. C. F. (Eds). 2005. Educação ambiental e sustentabilidade. Barueri SP: Manole. 878p. (Coleção Ambiental, 3)."
